$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3 currently holds "None"; update it to the new status value "TEST02-PC".
$ws.Range("D3").Value = "TEST02-PC"
